# Update values in column E ("Name of Algo" output values) to reflect
# new imputation results from the KNN algorithm run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value  = 17.346
$ws.Range("E18").Value = 16.593
$ws.Range("E20").Value = 16.291
$ws.Range("E27").Value = 16.531
$ws.Range("E69").Value = 17.4
$ws.Range("E76").Value = 16.558
$ws.Range("E82").Value = 16.776
